$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
for ($r = 2; $r -le 18; $r++) { $wsExpo.Cells.Item($r, 2).NumberFormat = "@" }
# Row 2: 合肥·CW国潮动漫游戏嘉年华
$wsExpo.Cells.Item(2, 1).Value = 1
$wsExpo.Cells.Item(2, 2).Value = '2024-03-16'
$wsExpo.Cells.Item(2, 3).Value = '合肥·CW国潮动漫游戏嘉年华'
$wsExpo.Cells.Item(2, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsExpo.Cells.Item(2, 5).Value = '2024.03.16 09:30-03.17 17:00'
$wsExpo.Cells.Item(2, 6).Value = 3487
$wsExpo.Cells.Item(2, 7).Value = 65
$wsExpo.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81284'
$wsExpo.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg'

# Row 3: 合肥·CW国潮动漫游戏嘉年华-赵路内场
$wsExpo.Cells.Item(3, 1).Value = 2
$wsExpo.Cells.Item(3, 2).Value = '2024-03-17'
$wsExpo.Cells.Item(3, 3).Value = '合肥·CW国潮动漫游戏嘉年华-赵路内场'
$wsExpo.Cells.Item(3, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsExpo.Cells.Item(3, 5).Value = '2024.03.17 09:00-03.17 17:00'
$wsExpo.Cells.Item(3, 6).Value = 746
$wsExpo.Cells.Item(3, 7).Value = '已售罄'
$wsExpo.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81954'
$wsExpo.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg'

# Row 4: 合肥·原&铁&崩 only展
$wsExpo.Cells.Item(4, 1).Value = 3
$wsExpo.Cells.Item(4, 2).Value = '2024-03-23'
$wsExpo.Cells.Item(4, 3).Value = '合肥·原&铁&崩 only展'
$wsExpo.Cells.Item(4, 4).Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsExpo.Cells.Item(4, 5).Value = '2024.03.23 09:00-03.23 17:00'
$wsExpo.Cells.Item(4, 6).Value = 142
$wsExpo.Cells.Item(4, 7).Value = 58
$wsExpo.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81574'
$wsExpo.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png'

# Row 5: 合肥· 第二届漫画城市动漫展 -故事再次开始
$wsExpo.Cells.Item(5, 1).Value = 4
$wsExpo.Cells.Item(5, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(5, 3).Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
$wsExpo.Cells.Item(5, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsExpo.Cells.Item(5, 5).Value = '2024.04.04 09:00-04.05 17:00'
$wsExpo.Cells.Item(5, 6).Value = 7009
$wsExpo.Cells.Item(5, 7).Value = 60
$wsExpo.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
$wsExpo.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg'

# Row 6: 合肥·环形宇宙动漫游戏嘉年华
$wsExpo.Cells.Item(6, 1).Value = 5
$wsExpo.Cells.Item(6, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(6, 3).Value = '合肥·环形宇宙动漫游戏嘉年华'
$wsExpo.Cells.Item(6, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsExpo.Cells.Item(6, 5).Value = '2024.04.04 09:30-04.05 17:00'
$wsExpo.Cells.Item(6, 6).Value = 2691
$wsExpo.Cells.Item(6, 7).Value = 65
$wsExpo.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81916'
$wsExpo.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg'

# Row 7: 合肥·环形宇宙动漫游戏嘉年华内场票-谢莹
$wsExpo.Cells.Item(7, 1).Value = 6
$wsExpo.Cells.Item(7, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(7, 3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-谢莹'
$wsExpo.Cells.Item(7, 4).Value = '锦绣大道3899号 合肥滨湖会展中心'
$wsExpo.Cells.Item(7, 5).Value = '2024.04.04 09:00-04.04 17:00'
$wsExpo.Cells.Item(7, 6).Value = 50
$wsExpo.Cells.Item(7, 7).Value = 118
$wsExpo.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82515'
$wsExpo.Cells.Item(7, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/L2DFEeao1709800386283.jpeg'

# Row 8: 合肥·环形宇宙动漫游戏嘉年华内场票-钱文青
$wsExpo.Cells.Item(8, 1).Value = 7
$wsExpo.Cells.Item(8, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(8, 3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-钱文青'
$wsExpo.Cells.Item(8, 4).Value = '锦绣大道3899号 合肥滨湖会展中心'
$wsExpo.Cells.Item(8, 5).Value = '2024.04.04 09:00-04.04 17:00'
$wsExpo.Cells.Item(8, 6).Value = 121
$wsExpo.Cells.Item(8, 7).Value = 238
$wsExpo.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82375'
$wsExpo.Cells.Item(8, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/3cHtIycW1709692273438.jpeg'

# Row 9: 合肥·第二届漫画城市动漫展内场-柯暮卿
$wsExpo.Cells.Item(9, 1).Value = 8
$wsExpo.Cells.Item(9, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(9, 3).Value = '合肥·第二届漫画城市动漫展内场-柯暮卿'
$wsExpo.Cells.Item(9, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsExpo.Cells.Item(9, 5).Value = '2024.04.04 10:00-04.04 17:00'
$wsExpo.Cells.Item(9, 6).Value = 29
$wsExpo.Cells.Item(9, 7).Value = 158
$wsExpo.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82192'
$wsExpo.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg'

# Row 10: 合肥·第二届漫画城市动漫展内场-风袖
$wsExpo.Cells.Item(10, 1).Value = 9
$wsExpo.Cells.Item(10, 2).Value = '2024-04-04'
$wsExpo.Cells.Item(10, 3).Value = '合肥·第二届漫画城市动漫展内场-风袖'
$wsExpo.Cells.Item(10, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsExpo.Cells.Item(10, 5).Value = '2024.04.04 10:00-04.04 17:00'
$wsExpo.Cells.Item(10, 6).Value = 38
$wsExpo.Cells.Item(10, 7).Value = 158
$wsExpo.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82191'
$wsExpo.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg'

# Row 11: 合肥· 第二届漫画城市动漫展内场-《琅声雅集》
$wsExpo.Cells.Item(11, 1).Value = 10
$wsExpo.Cells.Item(11, 2).Value = '2024-04-05'
$wsExpo.Cells.Item(11, 3).Value = '合肥· 第二届漫画城市动漫展内场-《琅声雅集》'
$wsExpo.Cells.Item(11, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsExpo.Cells.Item(11, 5).Value = '2024.04.05 10:00-04.05 17:00'
$wsExpo.Cells.Item(11, 6).Value = 85
$wsExpo.Cells.Item(11, 7).Value = 168
$wsExpo.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82189'
$wsExpo.Cells.Item(11, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/ZKkx4hTN1709191842946.jpeg'

# Row 12: 合肥·AOO动漫游戏嘉年华
$wsExpo.Cells.Item(12, 1).Value = 11
$wsExpo.Cells.Item(12, 2).Value = '2024-04-13'
$wsExpo.Cells.Item(12, 3).Value = '合肥·AOO动漫游戏嘉年华'
$wsExpo.Cells.Item(12, 4).Value = '芙蓉路291号 正大广场'
$wsExpo.Cells.Item(12, 5).Value = '2024.04.13 10:00-04.14 17:00'
$wsExpo.Cells.Item(12, 6).Value = 41
$wsExpo.Cells.Item(12, 7).Value = 49.9
$wsExpo.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82500'
$wsExpo.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/IpaidRy21709797005042.png'

# Row 13: 合肥·首届运动番only
$wsExpo.Cells.Item(13, 1).Value = 12
$wsExpo.Cells.Item(13, 2).Value = '2024-04-20'
$wsExpo.Cells.Item(13, 3).Value = '合肥·首届运动番only'
$wsExpo.Cells.Item(13, 4).Value = '繁华大道6177号 YONEX百胜羽毛球馆(包河店)'
$wsExpo.Cells.Item(13, 5).Value = '2024.04.20 10:00-04.20 17:00'
$wsExpo.Cells.Item(13, 6).Value = 1
$wsExpo.Cells.Item(13, 7).Value = 58
$wsExpo.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82924'
$wsExpo.Cells.Item(13, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/Vw8yXcUF1710489586295.jpeg'

# Row 14: 合肥·银魂only
$wsExpo.Cells.Item(14, 1).Value = 13
$wsExpo.Cells.Item(14, 2).Value = '2024-04-21'
$wsExpo.Cells.Item(14, 3).Value = '合肥·银魂only'
$wsExpo.Cells.Item(14, 4).Value = '濉溪路118号 合肥栢景假日酒店'
$wsExpo.Cells.Item(14, 5).Value = '2024.04.21 09:00-04.21 17:00'
$wsExpo.Cells.Item(14, 6).Value = 176
$wsExpo.Cells.Item(14, 7).Value = 55
$wsExpo.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82145'
$wsExpo.Cells.Item(14, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg'

# Row 15: 合肥·Look Look动漫嘉年华
$wsExpo.Cells.Item(15, 1).Value = 14
$wsExpo.Cells.Item(15, 2).Value = '2024-05-01'
$wsExpo.Cells.Item(15, 3).Value = '合肥·Look Look动漫嘉年华'
$wsExpo.Cells.Item(15, 4).Value = '新站区东方大道288号 少荃体育中心'
$wsExpo.Cells.Item(15, 5).Value = '2024.05.01 10:00-05.01 17:30'
$wsExpo.Cells.Item(15, 6).Value = 585
$wsExpo.Cells.Item(15, 7).Value = 29.9
$wsExpo.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82311'
$wsExpo.Cells.Item(15, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png'

# Row 16: 合肥·第十三届次元之门动漫游戏博览会
$wsExpo.Cells.Item(16, 1).Value = 15
$wsExpo.Cells.Item(16, 2).Value = '2024-05-01'
$wsExpo.Cells.Item(16, 3).Value = '合肥·第十三届次元之门动漫游戏博览会'
$wsExpo.Cells.Item(16, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsExpo.Cells.Item(16, 5).Value = '2024.05.01 10:00-05.03 17:00'
$wsExpo.Cells.Item(16, 6).Value = 1
$wsExpo.Cells.Item(16, 7).Value = '不可售'
$wsExpo.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82910'
$wsExpo.Cells.Item(16, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/NiDA76Th1710471014688.jpeg'

# Row 17: 合肥·BH动漫游戏展
$wsExpo.Cells.Item(17, 1).Value = 16
$wsExpo.Cells.Item(2, 1).Copy() | Out-Null
$wsExpo.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsExpo.Cells.Item(17, 1).Value = 16
$wsExpo.Cells.Item(17, 2).Value = '2024-05-03'
$wsExpo.Cells.Item(17, 3).Value = '合肥·BH动漫游戏展'
$wsExpo.Cells.Item(17, 4).Value = '科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地'
$wsExpo.Cells.Item(17, 5).Value = '2024.05.03 10:00-05.04 16:00'
$wsExpo.Cells.Item(17, 6).Value = 44
$wsExpo.Cells.Item(17, 7).Value = 40
$wsExpo.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82199'
$wsExpo.Cells.Item(17, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg'

# Row 18: 合肥·梦时空SPO1动漫展（取消）
$wsExpo.Cells.Item(18, 1).Value = 17
$wsExpo.Cells.Item(2, 1).Copy() | Out-Null
$wsExpo.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsExpo.Cells.Item(18, 1).Value = 17
$wsExpo.Cells.Item(18, 2).Value = '2024-05-18'
$wsExpo.Cells.Item(18, 3).Value = '合肥·梦时空SPO1动漫展（取消）'
$wsExpo.Cells.Item(18, 4).Value = '阜阳路16号 银瑞林国际大酒店'
$wsExpo.Cells.Item(18, 5).Value = '2024.05.18 10:00-05.18 17:00'
$wsExpo.Cells.Item(18, 6).Value = 131
$wsExpo.Cells.Item(18, 7).Value = '不可售'
$wsExpo.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
$wsExpo.Cells.Item(18, 9).Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

# --- Sheet "演出" (Show) ---
$wsShow = $wb.Worksheets.Item("演出")
for ($r = 2; $r -le 2; $r++) { $wsShow.Cells.Item($r, 2).NumberFormat = "@" }
# Row 2: 合肥·全国地下偶像联合公演展-永乐大典Vol.01
$wsShow.Cells.Item(2, 1).Value = 1
$wsShow.Cells.Item(2, 2).Value = '2024-03-16'
$wsShow.Cells.Item(2, 3).Value = '合肥·全国地下偶像联合公演展-永乐大典Vol.01'
$wsShow.Cells.Item(2, 4).Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsShow.Cells.Item(2, 5).Value = '2024.03.16 14:00-03.16 23:00'
$wsShow.Cells.Item(2, 6).Value = 26
$wsShow.Cells.Item(2, 7).Value = 78
$wsShow.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82442'
$wsShow.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/zi1Jk4QG1709716762992.jpeg'

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
for ($r = 2; $r -le 19; $r++) { $wsAll.Cells.Item($r, 2).NumberFormat = "@" }
# Row 2: 合肥·CW国潮动漫游戏嘉年华
$wsAll.Cells.Item(2, 1).Value = 1
$wsAll.Cells.Item(2, 2).Value = '2024-03-16'
$wsAll.Cells.Item(2, 3).Value = '合肥·CW国潮动漫游戏嘉年华'
$wsAll.Cells.Item(2, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsAll.Cells.Item(2, 5).Value = '2024.03.16 09:30-03.17 17:00'
$wsAll.Cells.Item(2, 6).Value = 3487
$wsAll.Cells.Item(2, 7).Value = 65
$wsAll.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81284'
$wsAll.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg'

# Row 3: 合肥·全国地下偶像联合公演展-永乐大典Vol.01
$wsAll.Cells.Item(3, 1).Value = 2
$wsAll.Cells.Item(3, 2).Value = '2024-03-16'
$wsAll.Cells.Item(3, 3).Value = '合肥·全国地下偶像联合公演展-永乐大典Vol.01'
$wsAll.Cells.Item(3, 4).Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsAll.Cells.Item(3, 5).Value = '2024.03.16 14:00-03.16 23:00'
$wsAll.Cells.Item(3, 6).Value = 26
$wsAll.Cells.Item(3, 7).Value = 78
$wsAll.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82442'
$wsAll.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/zi1Jk4QG1709716762992.jpeg'

# Row 4: 合肥·CW国潮动漫游戏嘉年华-赵路内场
$wsAll.Cells.Item(4, 1).Value = 3
$wsAll.Cells.Item(4, 2).Value = '2024-03-17'
$wsAll.Cells.Item(4, 3).Value = '合肥·CW国潮动漫游戏嘉年华-赵路内场'
$wsAll.Cells.Item(4, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsAll.Cells.Item(4, 5).Value = '2024.03.17 09:00-03.17 17:00'
$wsAll.Cells.Item(4, 6).Value = 746
$wsAll.Cells.Item(4, 7).Value = '已售罄'
$wsAll.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81954'
$wsAll.Cells.Item(4, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg'

# Row 5: 合肥·原&铁&崩 only展
$wsAll.Cells.Item(5, 1).Value = 4
$wsAll.Cells.Item(5, 2).Value = '2024-03-23'
$wsAll.Cells.Item(5, 3).Value = '合肥·原&铁&崩 only展'
$wsAll.Cells.Item(5, 4).Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsAll.Cells.Item(5, 5).Value = '2024.03.23 09:00-03.23 17:00'
$wsAll.Cells.Item(5, 6).Value = 142
$wsAll.Cells.Item(5, 7).Value = 58
$wsAll.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81574'
$wsAll.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png'

# Row 6: 合肥· 第二届漫画城市动漫展 -故事再次开始
$wsAll.Cells.Item(6, 1).Value = 5
$wsAll.Cells.Item(6, 2).Value = '2024-04-04'
$wsAll.Cells.Item(6, 3).Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
$wsAll.Cells.Item(6, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsAll.Cells.Item(6, 5).Value = '2024.04.04 09:00-04.05 17:00'
$wsAll.Cells.Item(6, 6).Value = 7009
$wsAll.Cells.Item(6, 7).Value = 60
$wsAll.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
$wsAll.Cells.Item(6, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg'

# Row 7: 合肥·环形宇宙动漫游戏嘉年华
$wsAll.Cells.Item(7, 1).Value = 6
$wsAll.Cells.Item(7, 2).Value = '2024-04-04'
$wsAll.Cells.Item(7, 3).Value = '合肥·环形宇宙动漫游戏嘉年华'
$wsAll.Cells.Item(7, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsAll.Cells.Item(7, 5).Value = '2024.04.04 09:30-04.05 17:00'
$wsAll.Cells.Item(7, 6).Value = 2692
$wsAll.Cells.Item(7, 7).Value = 65
$wsAll.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81916'
$wsAll.Cells.Item(7, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg'

# Row 8: 合肥·环形宇宙动漫游戏嘉年华内场票-谢莹
$wsAll.Cells.Item(8, 1).Value = 7
$wsAll.Cells.Item(8, 2).Value = '2024-04-04'
$wsAll.Cells.Item(8, 3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-谢莹'
$wsAll.Cells.Item(8, 4).Value = '锦绣大道3899号 合肥滨湖会展中心'
$wsAll.Cells.Item(8, 5).Value = '2024.04.04 09:00-04.04 17:00'
$wsAll.Cells.Item(8, 6).Value = 50
$wsAll.Cells.Item(8, 7).Value = 118
$wsAll.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82515'
$wsAll.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/L2DFEeao1709800386283.jpeg'

# Row 9: 合肥·环形宇宙动漫游戏嘉年华内场票-钱文青
$wsAll.Cells.Item(9, 1).Value = 8
$wsAll.Cells.Item(9, 2).Value = '2024-04-04'
$wsAll.Cells.Item(9, 3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-钱文青'
$wsAll.Cells.Item(9, 4).Value = '锦绣大道3899号 合肥滨湖会展中心'
$wsAll.Cells.Item(9, 5).Value = '2024.04.04 09:00-04.04 17:00'
$wsAll.Cells.Item(9, 6).Value = 121
$wsAll.Cells.Item(9, 7).Value = 238
$wsAll.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82375'
$wsAll.Cells.Item(9, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/3cHtIycW1709692273438.jpeg'

# Row 10: 合肥·第二届漫画城市动漫展内场-柯暮卿
$wsAll.Cells.Item(10, 1).Value = 9
$wsAll.Cells.Item(10, 2).Value = '2024-04-04'
$wsAll.Cells.Item(10, 3).Value = '合肥·第二届漫画城市动漫展内场-柯暮卿'
$wsAll.Cells.Item(10, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsAll.Cells.Item(10, 5).Value = '2024.04.04 10:00-04.04 17:00'
$wsAll.Cells.Item(10, 6).Value = 29
$wsAll.Cells.Item(10, 7).Value = 158
$wsAll.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82192'
$wsAll.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg'

# Row 11: 合肥·第二届漫画城市动漫展内场-风袖
$wsAll.Cells.Item(11, 1).Value = 10
$wsAll.Cells.Item(11, 2).Value = '2024-04-04'
$wsAll.Cells.Item(11, 3).Value = '合肥·第二届漫画城市动漫展内场-风袖'
$wsAll.Cells.Item(11, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsAll.Cells.Item(11, 5).Value = '2024.04.04 10:00-04.04 17:00'
$wsAll.Cells.Item(11, 6).Value = 38
$wsAll.Cells.Item(11, 7).Value = 158
$wsAll.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82191'
$wsAll.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg'

# Row 12: 合肥· 第二届漫画城市动漫展内场-《琅声雅集》
$wsAll.Cells.Item(12, 1).Value = 11
$wsAll.Cells.Item(12, 2).Value = '2024-04-05'
$wsAll.Cells.Item(12, 3).Value = '合肥· 第二届漫画城市动漫展内场-《琅声雅集》'
$wsAll.Cells.Item(12, 4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$wsAll.Cells.Item(12, 5).Value = '2024.04.05 10:00-04.05 17:00'
$wsAll.Cells.Item(12, 6).Value = 85
$wsAll.Cells.Item(12, 7).Value = 168
$wsAll.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82189'
$wsAll.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/ZKkx4hTN1709191842946.jpeg'

# Row 13: 合肥·AOO动漫游戏嘉年华
$wsAll.Cells.Item(13, 1).Value = 12
$wsAll.Cells.Item(13, 2).Value = '2024-04-13'
$wsAll.Cells.Item(13, 3).Value = '合肥·AOO动漫游戏嘉年华'
$wsAll.Cells.Item(13, 4).Value = '芙蓉路291号 正大广场'
$wsAll.Cells.Item(13, 5).Value = '2024.04.13 10:00-04.14 17:00'
$wsAll.Cells.Item(13, 6).Value = 41
$wsAll.Cells.Item(13, 7).Value = 49.9
$wsAll.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82500'
$wsAll.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/IpaidRy21709797005042.png'

# Row 14: 合肥·首届运动番only
$wsAll.Cells.Item(14, 1).Value = 13
$wsAll.Cells.Item(14, 2).Value = '2024-04-20'
$wsAll.Cells.Item(14, 3).Value = '合肥·首届运动番only'
$wsAll.Cells.Item(14, 4).Value = '繁华大道6177号 YONEX百胜羽毛球馆(包河店)'
$wsAll.Cells.Item(14, 5).Value = '2024.04.20 10:00-04.20 17:00'
$wsAll.Cells.Item(14, 6).Value = 1
$wsAll.Cells.Item(14, 7).Value = 58
$wsAll.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82924'
$wsAll.Cells.Item(14, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/Vw8yXcUF1710489586295.jpeg'

# Row 15: 合肥·银魂only
$wsAll.Cells.Item(15, 1).Value = 14
$wsAll.Cells.Item(15, 2).Value = '2024-04-21'
$wsAll.Cells.Item(15, 3).Value = '合肥·银魂only'
$wsAll.Cells.Item(15, 4).Value = '濉溪路118号 合肥栢景假日酒店'
$wsAll.Cells.Item(15, 5).Value = '2024.04.21 09:00-04.21 17:00'
$wsAll.Cells.Item(15, 6).Value = 176
$wsAll.Cells.Item(15, 7).Value = 55
$wsAll.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82145'
$wsAll.Cells.Item(15, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg'

# Row 16: 合肥·Look Look动漫嘉年华
$wsAll.Cells.Item(16, 1).Value = 15
$wsAll.Cells.Item(16, 2).Value = '2024-05-01'
$wsAll.Cells.Item(16, 3).Value = '合肥·Look Look动漫嘉年华'
$wsAll.Cells.Item(16, 4).Value = '新站区东方大道288号 少荃体育中心'
$wsAll.Cells.Item(16, 5).Value = '2024.05.01 10:00-05.01 17:30'
$wsAll.Cells.Item(16, 6).Value = 585
$wsAll.Cells.Item(16, 7).Value = 29.9
$wsAll.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82311'
$wsAll.Cells.Item(16, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png'

# Row 17: 合肥·第十三届次元之门动漫游戏博览会
$wsAll.Cells.Item(17, 1).Value = 16
$wsAll.Cells.Item(17, 2).Value = '2024-05-01'
$wsAll.Cells.Item(17, 3).Value = '合肥·第十三届次元之门动漫游戏博览会'
$wsAll.Cells.Item(17, 4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$wsAll.Cells.Item(17, 5).Value = '2024.05.01 10:00-05.03 17:00'
$wsAll.Cells.Item(17, 6).Value = 1
$wsAll.Cells.Item(17, 7).Value = '不可售'
$wsAll.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82910'
$wsAll.Cells.Item(17, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/NiDA76Th1710471014688.jpeg'

# Row 18: 合肥·BH动漫游戏展
$wsAll.Cells.Item(18, 1).Value = 17
$wsAll.Cells.Item(2, 1).Copy() | Out-Null
$wsAll.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsAll.Cells.Item(18, 1).Value = 17
$wsAll.Cells.Item(18, 2).Value = '2024-05-03'
$wsAll.Cells.Item(18, 3).Value = '合肥·BH动漫游戏展'
$wsAll.Cells.Item(18, 4).Value = '科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地'
$wsAll.Cells.Item(18, 5).Value = '2024.05.03 10:00-05.04 16:00'
$wsAll.Cells.Item(18, 6).Value = 44
$wsAll.Cells.Item(18, 7).Value = 40
$wsAll.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82199'
$wsAll.Cells.Item(18, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg'

# Row 19: 合肥·梦时空SPO1动漫展（取消）
$wsAll.Cells.Item(19, 1).Value = 18
$wsAll.Cells.Item(2, 1).Copy() | Out-Null
$wsAll.Cells.Item(19, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsAll.Cells.Item(19, 1).Value = 18
$wsAll.Cells.Item(19, 2).Value = '2024-05-18'
$wsAll.Cells.Item(19, 3).Value = '合肥·梦时空SPO1动漫展（取消）'
$wsAll.Cells.Item(19, 4).Value = '阜阳路16号 银瑞林国际大酒店'
$wsAll.Cells.Item(19, 5).Value = '2024.05.18 10:00-05.18 17:00'
$wsAll.Cells.Item(19, 6).Value = 131
$wsAll.Cells.Item(19, 7).Value = '不可售'
$wsAll.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
$wsAll.Cells.Item(19, 9).Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'
